$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.143.75'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '1.814.41'
$ws.Range('E3').Value = '  -1.64%  '
$ws.Range('E4').Value = '  +0.80%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.03'
$ws.Range('E5').Value = '  +2.06%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.91'
$ws.Range('E8').Value = '  -5.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.323'
$ws.Range('E9').Value = '  +5.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0684'
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').Value = '2.077.92'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = '1.809.14'
$ws.Range('E13').Value = '  -1.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.04'
$ws.Range('E14').Value = '  -4.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.659'
$ws.Range('E15').Value = '  -0.10%  '
$ws.Range('E16').Value = '  -1.92%  '
$ws.Range('D17').Value = '35.099.64'
$ws.Range('E17').Value = '  -0.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.56'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').Value = '0.0₃0790'
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '238.75'
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.88'
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.66'
$ws.Range('E22').Value = '  -1.01%  '
$ws.Range('E23').Value = '  +0.83%  '
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.72'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.84'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.48'
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('E29').Value = '  +20.03%  '
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('E31').Value = '  +5.42%  '
$ws.Range('D32').Value = '3.329.39'
$ws.Range('E32').Value = '  -8.69%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0553'
$ws.Range('E33').Value = '  +2.74%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.99'
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('E35').Value = '  -6.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.14'
$ws.Range('E36').Value = '  +5.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '92.40'
$ws.Range('E37').Value = '  +2.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.677'
$ws.Range('E38').Value = '  +0.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0193'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.312.37'
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.28'
$ws.Range('E41').Value = '  +2.00%  '
$ws.Range('E42').Value = '  -2.24%  '
$ws.Range('E43').Value = '  +1.32%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.49'
$ws.Range('E44').Value = '  -2.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.30'
$ws.Range('E45').Value = '  -5.25%  '
$ws.Range('E46').Value = '  -2.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.33'
$ws.Range('E47').Value = '  +4.39%  '
$ws.Range('E48').Value = '  -1.23%  '
$ws.Range('D49').Value = '1.992.68'
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('E50').Value = '  +0.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0655'
$ws.Range('E51').Value = '  +5.93%  '
